$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.588.50'
$ws.Range('E2').Value = '  -7.30%  '
$ws.Range('D3').Value = '1.693.49'
$ws.Range('E3').Value = '  -5.81%  '
$ws.Range('E4').Value = '  +0.34%  '
$ws.Range('D5').Value = '219.71'
$ws.Range('E5').Value = '  -5.17%  '
$ws.Range('D6').Value = '0.5109'
$ws.Range('E6').Value = '  -13.01%  '
$ws.Range('D7').Value = '1.007'
$ws.Range('E7').Value = '  +0.24%  '
$ws.Range('D8').Value = '0.2653'
$ws.Range('E8').Value = '  -4.21%  '
$ws.Range('D9').Value = '22.15'
$ws.Range('E9').Value = '  -4.58%  '
$ws.Range('D10').Value = '0.06312'
$ws.Range('E10').Value = '  -6.66%  '
$ws.Range('D11').Value = '0.07355'
$ws.Range('E11').Value = '  -2.24%  '
$ws.Range('D12').Value = '1.701.02'
$ws.Range('E12').Value = '  -5.22%  '
$ws.Range('D13').Value = '4.525'
$ws.Range('D14').Value = '0.5783'
$ws.Range('D15').Value = '1.926.28'
$ws.Range('E15').Value = '  -5.63%  '
$ws.Range('D16').Value = '0.000008483'
$ws.Range('E16').Value = '  -6.71%  '
$ws.Range('E17').Value = '  -13.02%  '
$ws.Range('D18').Value = '26.626.36'
$ws.Range('D19').Value = '4.982'
$ws.Range('E19').Value = '  -9.02%  '
$ws.Range('E20').Value = '  +0.22%  '
$ws.Range('D21').Value = '10.97'
$ws.Range('E21').Value = '  -4.61%  '
$ws.Range('E22').Value = '  -11.40%  '
$ws.Range('D23').Value = '6.254'
$ws.Range('E23').Value = '  -8.40%  '
$ws.Range('E24').Value = '  +0.31%  '
$ws.Range('D25').Value = '144.62'
$ws.Range('E25').Value = '  -5.83%  '
$ws.Range('D26').Value = '7.468'
$ws.Range('E26').Value = '  -7.60%  '
$ws.Range('D27').Value = '0.1163'
$ws.Range('E27').Value = '  -7.77%  '
$ws.Range('D28').Value = '15.89'
$ws.Range('E28').Value = '  -3.47%  '
$ws.Range('D29').Value = '1.341'
$ws.Range('E29').Value = '  -5.94%  '
$ws.Range('D30').Value = '0.05737'
$ws.Range('E30').Value = '  -6.28%  '
$ws.Range('D31').Value = '1.339'
$ws.Range('E31').Value = '  -5.88%  '
$ws.Range('D32').Value = '3.516'
$ws.Range('E32').Value = '  -7.30%  '
$ws.Range('D33').Value = '3.506'
$ws.Range('E33').Value = '  -7.99%  '
$ws.Range('D34').Value = '1.640'
$ws.Range('E34').Value = '  -5.64%  '
$ws.Range('E35').Value = '  -2.76%  '
$ws.Range('D36').Value = '0.5997'
$ws.Range('E36').Value = '  -6.60%  '
$ws.Range('E37').Value = '  -5.53%  '
$ws.Range('D38').Value = '2.686'
$ws.Range('E38').Value = '  -0.93%  '
$ws.Range('E39').Value = '  -4.67%  '
$ws.Range('D40').Value = '1.103.61'
$ws.Range('E40').Value = '  -3.39%  '
$ws.Range('D41').Value = '0.8569'
$ws.Range('E41').Value = '  -3.05%  '
$ws.Range('D42').Value = '5.832'
$ws.Range('E42').Value = '  -9.57%  '
$ws.Range('E43').Value = '  -0.34%  '
$ws.Range('D44').Value = '99.22'
$ws.Range('E44').Value = '  -0.88%  '
$ws.Range('D45').Value = '1.852.42'
$ws.Range('E45').Value = '  -5.09%  '
$ws.Range('D46').Value = '0.00000000118'
$ws.Range('E46').Value = '  +5.30%  '
$ws.Range('D47').Value = '56.61'
$ws.Range('E47').Value = '  -5.76%  '
$ws.Range('D48').Value = '1.005'
$ws.Range('E48').Value = '  +0.71%  '
$ws.Range('D49').Value = '8.130'
$ws.Range('E49').Value = '  -2.68%  '
$ws.Range('D50').Value = '0.4328'
$ws.Range('E50').Value = '  -3.42%  '
$ws.Range('D51').Value = '0.05234'
$ws.Range('E51').Value = '  -4.62%  '
